# ---------------------------------------------------------------------------
# Edit: add a new "Player Info" sheet (as the first sheet) and update the
# MATCH_CARD_LINK column (header + values) on the "ODI Batting" and
# "ODI Bowling" sheets to a new MATCH_CODE column containing just the
# numeric match code instead of the full scorecard URL.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# --- helper: write a value as TEXT (not auto-converted to a number) without
# leaving a stray style behind -----------------------------------------------
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

$battingCodes = @(
    @(2, "4628"),
    @(3, "4679"),
    @(4, "4682"),
    @(5, "4685"),
    @(6, "4717"),
    @(7, "4726"),
    @(8, "4729"),
    @(9, "4734")
)
foreach ($pair in $battingCodes) {
    Set-TextValue $batting.Range("D$($pair[0])") $pair[1]
}

# ---------------------------------------------------------------------------
# 2. Update "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE
# ---------------------------------------------------------------------------
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @(
    @(2, "4628"),
    @(3, "4679"),
    @(4, "4682"),
    @(5, "4685"),
    @(6, "4717"),
    @(7, "4726"),
    @(8, "4734")
)
foreach ($pair in $bowlingCodes) {
    Set-TextValue $bowling.Range("B$($pair[0])") $pair[1]
}

# ---------------------------------------------------------------------------
# 3. Insert a new "Player Info" sheet as the first sheet in the workbook
# ---------------------------------------------------------------------------
$info = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$info.Name = "Player Info"

$info.Range("A1").Value = "ID"
$info.Range("B1").Value = "NAME"
$info.Range("C1").Value = "BATTING_HAND"
$info.Range("D1").Value = "BOWL_STYLE"

$headerRange = $info.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Borders.LineStyle = 1
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160

Set-TextValue $info.Range("A2") "4925"
$info.Range("B2").Value = "Ebadot Hossain Chowdhury"
$info.Range("C2").Value = "Right Handed"
$info.Range("D2").Value = "Right Arm Medium"
